$wb = $excel.ActiveWorkbook

# Rename sheet "MMAP4P#1" -> "MMAP4P"
$mma = $wb.Worksheets.Item("MMAP4P#1")
$mma.Name = "MMAP4P"

# Update the MMAP4P sheet's content/selection/column width
$mma.Range("A2").Value = "MMA pound-for-pound rankings"
$mma.Columns.Item(1).ColumnWidth = 33.17
$mma.Range("A2").Select()

# Update the NHLScoresText sheet's content/selection (becomes the active sheet/tab)
$nhl = $wb.Worksheets.Item("NHLScoresText")
$nhl.Range("A2").Value = "NHL Scoreboard"
$nhl.Range("A2").Select()
